$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the corrected guest lecture date entry for row 14 (week 15)
$ws.Range("D14").Value = "15.04: Gjesteforelesning med Ole-Petter Hansen, Tryg Forsikring  på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."

# Update the active selection to reflect the last edited cell
$ws.Range("D14").Select()
